# Update Mappings 22 Ontologies
#
# - fixes the capitalisation of the metadata4ing_IRI / metadata4ing_DESC
#   headers
# - adds a new "VIMMP_DEF" column (F)
# - inserts a new "MeasurementUnit" (EMMO) row between "Agent" and
#   "Variable"
# - appends a new "NumericalAssignment" row at the end
# - refreshes every hyperlink so the rIds / locations line up with the new
#   row layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1).
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,2).Value = "metadata4ing_IRI"
$ws.Cells.Item(1,3).Value = "metadata4ing_DESC"
$ws.Cells.Item(1,6).Value = "VIMMP_DEF"

# ---------------------------------------------------------------------------
# 2. Data rows (2-9).
# ---------------------------------------------------------------------------
$data = @(
    @{ Row=2; A=0; B="http://purl.obolibrary.org/obo/BFO_0000015";
       C="{'label': 'process', 'prefLabel': 'process', 'altLabel': None, 'name': 'BFO_0000015'}";
       D="https://emmc.eu/semantics/evmpo/evmpo.ttl#process"; E="{'name': 'process'}"; F="[]" },
    @{ Row=3; A=1; B="http://xmlns.com/foaf/0.1/Person";
       C="{'label': 'Person', 'prefLabel': 'Person', 'altLabel': None, 'name': 'Person'}";
       D="https://purl.vimmp.eu/semantics/vico/vico.ttl#person"; E="{'name': 'Person'}"; F="[]" },
    @{ Row=4; A=2; B="http://xmlns.com/foaf/0.1/Agent";
       C="{'label': 'Agent', 'prefLabel': 'Agent', 'altLabel': None, 'name': 'Agent'}";
       D="https://emmc.eu/semantics/evmpo/evmpo.ttl#agent"; E="{'name': 'Agent'}"; F="[]" },
    @{ Row=5; A=3; B="http://emmo.info/emmo#EMMO_b081b346_7279_46ef_9a3d_2c088fcd79f4";
       C="{'label': None, 'prefLabel': 'MeasurementUnit', 'altLabel': None, 'name': 'EMMO_b081b346_7279_46ef_9a3d_2c088fcd79f4'}";
       D="https://purl.vimmp.eu/semantics/alignment/emmo1s.ttl#MeasurementUnit";
       E="{'label': 'MeasurementUnit', 'name': 'MeasurementUnit'}"; F="[]" },
    @{ Row=6; A=4; B="http://www.molmod.info/semantics/pims-ii.ttl#Variable";
       C="{'label': None, 'prefLabel': 'Variable', 'altLabel': None, 'name': 'Variable'}";
       D="https://purl.vimmp.eu/semantics/alignment/emmo1s.ttl#Variable";
       E="{'label': 'Variable', 'name': 'Variable'}"; F="[]" },
    @{ Row=7; A=5; B="https://schema.org/Project";
       C="{'label': None, 'prefLabel': 'Project', 'altLabel': None, 'name': 'Project'}";
       D="https://purl.vimmp.eu/semantics/mmto/mmto.ttl#project"; E="{'name': 'Project'}"; F="[]" },
    @{ Row=8; A=6; B="http://www.molmod.info/semantics/pims-ii.ttl#Property";
       C="{'label': None, 'prefLabel': 'Eigenschaft', 'altLabel': None, 'name': 'Property'}";
       D="https://emmc.eu/semantics/evmpo/evmpo.ttl#property"; E="{'name': 'Property'}"; F="[]" },
    @{ Row=9; A=7; B="http://w3id.org/nfdi4ing/metadata4ing#NumericalAssignment";
       C="{'label': None, 'prefLabel': 'numerical', 'altLabel': None, 'name': 'NumericalAssignment'}";
       D="https://purl.vimmp.eu/semantics/alignment/emmo1s.ttl#Numerical";
       E="{'label': 'numerical', 'name': 'numerical'}"; F="[]" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r,1).Value = $item.A
    $ws.Cells.Item($r,2).Value = $item.B
    $ws.Cells.Item($r,3).Value = $item.C
    $ws.Cells.Item($r,4).Value = $item.D
    $ws.Cells.Item($r,5).Value = $item.E
    $ws.Cells.Item($r,6).Value = $item.F
}

# ---------------------------------------------------------------------------
# 3. Hyperlinks. The engine does not shift hyperlinks when rows/values move,
#    so clear every hyperlink and recreate them all, in the exact order they
#    appear in the final sheet. That reproduces the expected rId1..rId16
#    relationship numbering.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Cells.Item(2,2), "http://purl.obolibrary.org/obo/BFO_0000015")
$ws.Hyperlinks.Add($ws.Cells.Item(2,4), "https://emmc.eu/semantics/evmpo/evmpo.ttl", "process")

$ws.Hyperlinks.Add($ws.Cells.Item(3,2), "http://xmlns.com/foaf/0.1/Person")
$ws.Hyperlinks.Add($ws.Cells.Item(3,4), "https://purl.vimmp.eu/semantics/vico/vico.ttl", "person")

$ws.Hyperlinks.Add($ws.Cells.Item(4,2), "http://xmlns.com/foaf/0.1/Agent")
$ws.Hyperlinks.Add($ws.Cells.Item(4,4), "https://emmc.eu/semantics/evmpo/evmpo.ttl", "agent")

$ws.Hyperlinks.Add($ws.Cells.Item(5,2), "http://emmo.info/emmo", "EMMO_b081b346_7279_46ef_9a3d_2c088fcd79f4")
$ws.Hyperlinks.Add($ws.Cells.Item(5,4), "https://purl.vimmp.eu/semantics/alignment/emmo1s.ttl", "MeasurementUnit")

$ws.Hyperlinks.Add($ws.Cells.Item(6,2), "http://www.molmod.info/semantics/pims-ii.ttl", "Variable")
$ws.Hyperlinks.Add($ws.Cells.Item(6,4), "https://purl.vimmp.eu/semantics/alignment/emmo1s.ttl", "Variable")

$ws.Hyperlinks.Add($ws.Cells.Item(7,2), "https://schema.org/Project")
$ws.Hyperlinks.Add($ws.Cells.Item(7,4), "https://purl.vimmp.eu/semantics/mmto/mmto.ttl", "project")

$ws.Hyperlinks.Add($ws.Cells.Item(8,2), "http://www.molmod.info/semantics/pims-ii.ttl", "Property")
$ws.Hyperlinks.Add($ws.Cells.Item(8,4), "https://emmc.eu/semantics/evmpo/evmpo.ttl", "property")

$ws.Hyperlinks.Add($ws.Cells.Item(9,2), "http://w3id.org/nfdi4ing/metadata4ing", "NumericalAssignment")
$ws.Hyperlinks.Add($ws.Cells.Item(9,4), "https://purl.vimmp.eu/semantics/alignment/emmo1s.ttl", "Numerical")

# ---------------------------------------------------------------------------
# 4. Re-apply the original cell styling. Adding a hyperlink re-derives the
#    cell's style, so every touched cell is re-stamped with the canonical
#    "index" style (column A), and "hyperlink" style (columns B & D) taken
#    from row 2, which keeps the whole table visually consistent and reuses
#    the existing style indexes instead of minting new ones.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,5).Copy()
$ws.Cells.Item(1,6).PasteSpecial(-4122)

foreach ($r in @(2,3,4,5,6,7,8,9)) {
    $ws.Cells.Item(2,1).Copy()
    $ws.Cells.Item($r,1).PasteSpecial(-4122)

    $ws.Cells.Item(2,2).Copy()
    $ws.Cells.Item($r,2).PasteSpecial(-4122)

    $ws.Cells.Item(2,4).Copy()
    $ws.Cells.Item($r,4).PasteSpecial(-4122)
}
